$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Color value (case change BLACK -> Black)
$ws.Range("B2").Value = "Black"

# Update Quantity value (3 -> 1), stored as text per the cell's text number format
$ws.Range("D2").Value = "1"

# Update the active selection to rows 3:38 (A3 active cell, A3:XFD38 selected)
$ws.Rows("3:38").Select() | Out-Null
